$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 65.93
$ws.Range("C5").Value = 7.76262235641479
$ws.Range("D5").Value = 1.19393062591553

$ws.Range("F5").Value = 64.66
$ws.Range("G5").Value = 6.74870824813843
$ws.Range("H5").Value = 1.30078768730164

$ws.Range("J5").Value = 64.66
$ws.Range("K5").Value = 6.74870824813843
$ws.Range("L5").Value = 1.30078768730164

$ws.Range("N5").Value = 65.06
$ws.Range("O5").Value = 7.98389196395874
$ws.Range("P5").Value = 2.12885189056396

$ws.Range("R5").Value = 65.06
$ws.Range("S5").Value = 7.98389196395874
$ws.Range("T5").Value = 2.12885189056396

$ws.Range("W5").Value = 7.42410612106323
$ws.Range("X5").Value = 2.03520464897156

$ws.Range("AA5").Value = 7.42410612106323
$ws.Range("AB5").Value = 2.03520464897156
